$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7314.5454
$ws.Range("J51").Value = 7243.1577
$ws.Range("L51").Value = 7243.1577
$ws.Range("N51").Value = -8211.1577
$ws.Range("H87").Value = 88999.5
$ws.Range("J87").Value = 88999.5
$ws.Range("L87").Value = 88999.5
$ws.Range("N87").Value = -91495.5
$ws.Range("H90").Value = 88999.5
$ws.Range("J90").Value = 88999.5
$ws.Range("L90").Value = 266998.5
$ws.Range("N90").Value = -279478.5
$ws.Range("H92").Value = 1238.1875
$ws.Range("I92").Value = 1434.5834
$ws.Range("K92").Value = 1434.5834
$ws.Range("M92").Value = -186.5834
$ws.Range("H101").Value = 41668732
$ws.Range("I101").Value = 55556020
$ws.Range("K101").Value = 166668060
$ws.Range("M101").Value = -166666438
$ws.Range("H112").Value = 22418.4
$ws.Range("J112").Value = 35369
$ws.Range("L112").Value = 106107
$ws.Range("N112").Value = -108323
$ws.Range("H127").Value = 1924.3
$ws.Range("I127").Value = 1392.875
$ws.Range("K127").Value = 4178.625
$ws.Range("M127").Value = 781.375
$ws.Range("H131").Value = 2817.963
$ws.Range("I131").Value = 2253.682
$ws.Range("K131").Value = 6761.045999999999
$ws.Range("M131").Value = -1721.045999999999
$ws.Range("H135").Value = 892.5
$ws.Range("I135").Value = 874.5172
$ws.Range("J135").Value = 967
$ws.Range("K135").Value = 7870.6548
$ws.Range("L135").Value = 8703
$ws.Range("M135").Value = -5335.6548
$ws.Range("N135").Value = -13773
$ws.Range("H138").Value = 3022.8298
$ws.Range("I138").Value = 1371.25
$ws.Range("K138").Value = 4113.75
$ws.Range("M138").Value = 1026.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1022.5
$ws.Range("I4").Value = 1296.2858
$ws.Range("K4").Value = 1296.2858
$ws.Range("M4").Value = -1180.2858
$ws.Range("H132").Value = 3318.8696
$ws.Range("I132").Value = 2332.2307
$ws.Range("J132").Value = 4601.5
$ws.Range("K132").Value = 6996.6921
$ws.Range("L132").Value = 13804.5
$ws.Range("M132").Value = -4466.6921
$ws.Range("N132").Value = -18864.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2759.5
$ws.Range("I20").Value = 1804.0769
$ws.Range("K20").Value = 1804.0769
$ws.Range("M20").Value = -1557.0769
$ws.Range("H86").Value = 4005209.2
$ws.Range("I86").Value = 5005946.5
$ws.Range("J86").Value = 2260.8
$ws.Range("K86").Value = 5005946.5
$ws.Range("L86").Value = 2260.8
$ws.Range("M86").Value = -5004823.5
$ws.Range("N86").Value = -4506.8
$ws.Range("H89").Value = 4005209.2
$ws.Range("I89").Value = 5005946.5
$ws.Range("J89").Value = 2260.8
$ws.Range("K89").Value = 25029732.5
$ws.Range("L89").Value = 11304
$ws.Range("M89").Value = -25024116.5
$ws.Range("N89").Value = -22536
$ws.Range("H94").Value = 3966024
$ws.Range("I94").Value = 6994222
$ws.Range("J94").Value = 29366.9
$ws.Range("K94").Value = 6994222
$ws.Range("L94").Value = 29366.9
$ws.Range("M94").Value = -6993771
$ws.Range("N94").Value = -30268.9
$ws.Range("H105").Value = 6946142
$ws.Range("I105").Value = 7814160
$ws.Range("K105").Value = 7814160
$ws.Range("M105").Value = -7812413

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1309.4286
$ws.Range("J94").Value = 1359.5
$ws.Range("L94").Value = 1359.5
$ws.Range("N94").Value = -2261.5
$ws.Range("H105").Value = 1727.3158
$ws.Range("I105").Value = 1145.625
$ws.Range("K105").Value = 1145.625
$ws.Range("M105").Value = 601.375
$ws.Range("H132").Value = 65264.367
$ws.Range("I132").Value = 49859.383
$ws.Range("K132").Value = 149578.149
$ws.Range("M132").Value = -147048.149

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3720926.2
$ws.Range("I4").Value = 7374185
$ws.Range("J4").Value = 67667.25
$ws.Range("K4").Value = 22122555
$ws.Range("L4").Value = 203001.75
$ws.Range("M4").Value = -22122443
$ws.Range("N4").Value = -203225.75
$ws.Range("H12").Value = 59481.332
$ws.Range("J12").Value = 265.42856
$ws.Range("L12").Value = 796.28568
$ws.Range("N12").Value = -1142.28568
$ws.Range("H45").Value = 5691.125
$ws.Range("J45").Value = 6357
$ws.Range("L45").Value = 19071
$ws.Range("N45").Value = -20135
$ws.Range("H60").Value = 1245
$ws.Range("I60").Value = 495
$ws.Range("K60").Value = 1485
$ws.Range("M60").Value = -1234
$ws.Range("H125").Value = 6419.9
$ws.Range("I125").Value = 2687.5
$ws.Range("K125").Value = 8062.5
$ws.Range("M125").Value = -3142.5
$ws.Range("H131").Value = 20835594
$ws.Range("I131").Value = 9259722
$ws.Range("J131").Value = 30306764
$ws.Range("K131").Value = 27779166
$ws.Range("L131").Value = 90920292
$ws.Range("M131").Value = -27774126
$ws.Range("N131").Value = -90930372
$ws.Range("H132").Value = 1608.625
$ws.Range("I132").Value = 1204.1111
$ws.Range("J132").Value = 2128.7144
$ws.Range("K132").Value = 10836.9999
$ws.Range("L132").Value = 19158.4296
$ws.Range("M132").Value = -8306.999900000001
$ws.Range("N132").Value = -24218.4296

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6905
$ws.Range("I2").Value = 1147.8889
$ws.Range("K2").Value = 1147.8889
$ws.Range("M2").Value = -1034.8889
$ws.Range("H70").Value = 8338430.5
$ws.Range("I70").Value = 10531136
$ws.Range("J70").Value = 6151.8
$ws.Range("K70").Value = 10531136
$ws.Range("L70").Value = 6151.8
$ws.Range("M70").Value = -10530866
$ws.Range("N70").Value = -6691.8
$ws.Range("H73").Value = 8338430.5
$ws.Range("I73").Value = 10531136
$ws.Range("J73").Value = 6151.8
$ws.Range("K73").Value = 10531136
$ws.Range("L73").Value = 6151.8
$ws.Range("M73").Value = -10530200
$ws.Range("N73").Value = -8023.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6970.0586
$ws.Range("I46").Value = 6849.3
$ws.Range("J46").Value = 7142.5713
$ws.Range("K46").Value = 6849.3
$ws.Range("L46").Value = 7142.5713
$ws.Range("M46").Value = -6661.3
$ws.Range("N46").Value = -7518.5713
$ws.Range("H48").Value = 35000
$ws.Range("J48").Value = 35000
$ws.Range("L48").Value = 35000
$ws.Range("N48").Value = -36322
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H55").Value = 1434.919
$ws.Range("I55").Value = 1421.5555
$ws.Range("K55").Value = 1421.5555
$ws.Range("M55").Value = -1248.5555
$ws.Range("H61").Value = 55557210
$ws.Range("I61").Value = 111111110
$ws.Range("J61").Value = 3305
$ws.Range("K61").Value = 111111110
$ws.Range("L61").Value = 3305
$ws.Range("M61").Value = -111110908
$ws.Range("N61").Value = -3709
$ws.Range("H113").Value = 55557210
$ws.Range("I113").Value = 111111110
$ws.Range("J113").Value = 3305
$ws.Range("K113").Value = 111111110
$ws.Range("L113").Value = 3305
$ws.Range("M113").Value = -111108940
$ws.Range("N113").Value = -7645

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2999.5
$ws.Range("J96").Value = 2999.5
$ws.Range("L96").Value = 2999.5
$ws.Range("N96").Value = -5745.5
$ws.Range("H113").Value = 1318.9
$ws.Range("I113").Value = 1073
$ws.Range("J113").Value = 1482.8334
$ws.Range("K113").Value = 3219
$ws.Range("L113").Value = 4448.5002
$ws.Range("M113").Value = -1049
$ws.Range("N113").Value = -8788.5002
$ws.Range("H122").Value = 1692.1154
$ws.Range("I122").Value = 1273
$ws.Range("K122").Value = 3819
$ws.Range("M122").Value = -1369
